$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 121.77778
$ws.Range("I5").Value = 121.25
$ws.Range("J5").Value = 126
$ws.Range("K5").Value = 121.25
$ws.Range("L5").Value = 126
$ws.Range("M5").Value = -6.25
$ws.Range("N5").Value = -356
$ws.Range("H64").Value = 4950
$ws.Range("J64").Value = 3200
$ws.Range("L64").Value = 3200
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 4950
$ws.Range("J67").Value = 3200
$ws.Range("L67").Value = 3200
$ws.Range("N67").Value = -4916
$ws.Range("H70").Value = 20000.666
$ws.Range("I70").Value = 9999.5
$ws.Range("K70").Value = 29998.5
$ws.Range("M70").Value = -29728.5
$ws.Range("H73").Value = 20000.666
$ws.Range("I73").Value = 9999.5
$ws.Range("K73").Value = 29998.5
$ws.Range("M73").Value = -29062.5
$ws.Range("H76").Value = 7122.1113
$ws.Range("I76").Value = 5099.5
$ws.Range("K76").Value = 5099.5
$ws.Range("M76").Value = -4784.5
$ws.Range("H79").Value = 7122.1113
$ws.Range("I79").Value = 5099.5
$ws.Range("K79").Value = 5099.5
$ws.Range("M79").Value = -4007.5
$ws.Range("H92").Value = 648.43475
$ws.Range("I92").Value = 553.58826
$ws.Range("K92").Value = 553.58826
$ws.Range("M92").Value = 694.41174
$ws.Range("H96").Value = 4017.25
$ws.Range("I96").Value = 3866.6667
$ws.Range("J96").Value = 4107.6
$ws.Range("K96").Value = 11600.0001
$ws.Range("L96").Value = 12322.8
$ws.Range("M96").Value = -10227.0001
$ws.Range("N96").Value = -15068.8
$ws.Range("H118").Value = 204
$ws.Range("I118").Value = 204
$ws.Range("K118").Value = 612
$ws.Range("M118").Value = 1045
$ws.Range("H141").Value = 9585.700000000001
$ws.Range("I141").Value = 9317.444
$ws.Range("K141").Value = 27952.332
$ws.Range("M141").Value = -22772.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10548.9795
$ws.Range("I32").Value = 9202.412
$ws.Range("J32").Value = 13819.214
$ws.Range("K32").Value = 9202.412
$ws.Range("L32").Value = 13819.214
$ws.Range("M32").Value = -8915.412
$ws.Range("N32").Value = -14393.214
$ws.Range("H55").Value = 69999
$ws.Range("J55").Value = 69999
$ws.Range("L55").Value = 69999
$ws.Range("N55").Value = -70629
$ws.Range("H97").Value = 4544.8335
$ws.Range("J97").Value = 11673.667
$ws.Range("L97").Value = 11673.667
$ws.Range("N97").Value = -12665.667
$ws.Range("H102").Value = 1285.1333
$ws.Range("I102").Value = 1283.25
$ws.Range("K102").Value = 1283.25
$ws.Range("M102").Value = 338.75
$ws.Range("H132").Value = 2464.2
$ws.Range("I132").Value = 2464.2
$ws.Range("K132").Value = 7392.599999999999
$ws.Range("M132").Value = -4862.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1132.4
$ws.Range("I94").Value = 1217.5294
$ws.Range("K94").Value = 1217.5294
$ws.Range("M94").Value = -766.5293999999999
$ws.Range("H99").Value = 3481.647
$ws.Range("J99").Value = 2828.8333
$ws.Range("L99").Value = 2828.8333
$ws.Range("N99").Value = -5824.8333
$ws.Range("H105").Value = 2740.9167
$ws.Range("I105").Value = 2083.0312
$ws.Range("K105").Value = 2083.0312
$ws.Range("M105").Value = -336.0311999999999
$ws.Range("H107").Value = 680.9048
$ws.Range("I107").Value = 682.35297
$ws.Range("J107").Value = 674.75
$ws.Range("K107").Value = 682.35297
$ws.Range("L107").Value = 674.75
$ws.Range("M107").Value = 1237.64703
$ws.Range("N107").Value = -4514.75
$ws.Range("H134").Value = 2938.1428
$ws.Range("I134").Value = 2101.9092
$ws.Range("K134").Value = 6305.7276
$ws.Range("M134").Value = -3770.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4870.077
$ws.Range("I31").Value = 3151.875
$ws.Range("J31").Value = 7619.2
$ws.Range("K31").Value = 3151.875
$ws.Range("L31").Value = 7619.2
$ws.Range("M31").Value = -2856.875
$ws.Range("N31").Value = -8209.200000000001
$ws.Range("H34").Value = 4870.077
$ws.Range("I34").Value = 3151.875
$ws.Range("J34").Value = 7619.2
$ws.Range("K34").Value = 3151.875
$ws.Range("L34").Value = 7619.2
$ws.Range("M34").Value = -2949.875
$ws.Range("N34").Value = -8023.2
$ws.Range("H58").Value = 5593.9287
$ws.Range("I58").Value = 4898
$ws.Range("J58").Value = 5980.5557
$ws.Range("K58").Value = 4898
$ws.Range("L58").Value = 5980.5557
$ws.Range("M58").Value = -4695
$ws.Range("N58").Value = -6386.5557
$ws.Range("H86").Value = 8974.5
$ws.Range("I86").Value = 7359.2
$ws.Range("J86").Value = 11666.667
$ws.Range("K86").Value = 7359.2
$ws.Range("L86").Value = 11666.667
$ws.Range("M86").Value = -6236.2
$ws.Range("N86").Value = -13912.667
$ws.Range("H89").Value = 8974.5
$ws.Range("I89").Value = 7359.2
$ws.Range("J89").Value = 11666.667
$ws.Range("K89").Value = 36796
$ws.Range("L89").Value = 58333.335
$ws.Range("M89").Value = -31180
$ws.Range("N89").Value = -69565.33499999999
$ws.Range("H105").Value = 454.27274
$ws.Range("I105").Value = 454.27274
$ws.Range("K105").Value = 454.27274
$ws.Range("M105").Value = 1292.72726
$ws.Range("H134").Value = 2341.2703
$ws.Range("I134").Value = 1788.2667
$ws.Range("K134").Value = 5364.800099999999
$ws.Range("M134").Value = -2829.800099999999
$ws.Range("H136").Value = 5593.9287
$ws.Range("I136").Value = 4898
$ws.Range("J136").Value = 5980.5557
$ws.Range("K136").Value = 14694
$ws.Range("L136").Value = 17941.6671
$ws.Range("M136").Value = -12144
$ws.Range("N136").Value = -23041.6671
$ws.Range("H141").Value = 139249
$ws.Range("J141").Value = 139249
$ws.Range("L141").Value = 139249
$ws.Range("N141").Value = -149609

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2001
$ws.Range("I80").Value = 2001
$ws.Range("K80").Value = 6003
$ws.Range("M80").Value = -5067
$ws.Range("H83").Value = 2001
$ws.Range("I83").Value = 2001
$ws.Range("K83").Value = 18009
$ws.Range("M83").Value = -13329
$ws.Range("H107").Value = 558.6896400000001
$ws.Range("J107").Value = 531.75
$ws.Range("L107").Value = 1595.25
$ws.Range("N107").Value = -5435.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 269.625
$ws.Range("I2").Value = 75.40000000000001
$ws.Range("J2").Value = 593.3333
$ws.Range("K2").Value = 75.40000000000001
$ws.Range("L2").Value = 593.3333
$ws.Range("M2").Value = 37.59999999999999
$ws.Range("N2").Value = -819.3333
$ws.Range("H93").Value = 69993.336
$ws.Range("J93").Value = 59990
$ws.Range("L93").Value = 59990
$ws.Range("N93").Value = -63734
$ws.Range("H132").Value = 1535.9556
$ws.Range("I132").Value = 1194.0233
$ws.Range("J132").Value = 8887.5
$ws.Range("K132").Value = 3582.0699
$ws.Range("L132").Value = 26662.5
$ws.Range("M132").Value = -1052.0699
$ws.Range("N132").Value = -31722.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4368.4
$ws.Range("I122").Value = 4460.625
$ws.Range("K122").Value = 13381.875
$ws.Range("M122").Value = -10931.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1497
$ws.Range("I96").Value = 1497.6666
$ws.Range("K96").Value = 1497.6666
$ws.Range("M96").Value = -124.6666
$ws.Range("H136").Value = 8990.799999999999
$ws.Range("I136").Value = 8989.75
$ws.Range("K136").Value = 26969.25
$ws.Range("M136").Value = -24419.25
